$wb = $excel.ActiveWorkbook

# --- CMT(throat) sheet: update header text and change the active selection ---
$wsThroat = $wb.Worksheets.Item("CMT(throat)")
$wsThroat.Range("C1").Value = "Manitenance edges"
$wsThroat.Range("B13").Select()

# --- PMT sheet becomes the active/selected tab (was CMT(platform)) ---
$wsPMT = $wb.Worksheets.Item("PMT")
$wsPMT.Activate()
